# Auto-generated edit script applying scheduled market-data refresh
# to the per-sheet Leve profit tables (columns H-N).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 3848
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3848
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 11544
$ws.Range("N46").Value = -11782
$ws.Range("H60").Value = 3848
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 3848
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 11544
$ws.Range("N60").Value = -12512
$ws.Range("H62").Value = 72982.47
$ws.Range("I62").Value = 107059.4
$ws.Range("J62").Value = 4828.6
$ws.Range("K62").Value = 107059.4
$ws.Range("L62").Value = 4828.6
$ws.Range("M62").Value = -106435.4
$ws.Range("N62").Value = -6076.6
$ws.Range("H65").Value = 72982.47
$ws.Range("I65").Value = 107059.4
$ws.Range("J65").Value = 4828.6
$ws.Range("K65").Value = 535297
$ws.Range("L65").Value = 24143
$ws.Range("M65").Value = -532177
$ws.Range("N65").Value = -30383
$ws.Range("H100").Value = 3413.2778
$ws.Range("I100").Value = 2981.5715
$ws.Range("J100").Value = 4924.25
$ws.Range("K100").Value = 2981.5715
$ws.Range("L100").Value = 4924.25
$ws.Range("M100").Value = -2440.5715
$ws.Range("N100").Value = -6006.25
$ws.Range("H127").Value = 98868.36
$ws.Range("I127").Value = 153405
$ws.Range("J127").Value = 702.4
$ws.Range("K127").Value = 460215
$ws.Range("L127").Value = 2107.2
$ws.Range("M127").Value = -455255
$ws.Range("N127").Value = -12027.2
$ws.Range("H132").Value = 1313.2391
$ws.Range("I132").Value = 893.0952
$ws.Range("J132").Value = 5724.75
$ws.Range("K132").Value = 2679.2856
$ws.Range("L132").Value = 17174.25
$ws.Range("M132").Value = -149.2856000000002
$ws.Range("N132").Value = -22234.25
$ws.Range("H133").Value = 80841.3
$ws.Range("I133").Value = 83709
$ws.Range("J133").Value = 80662.06
$ws.Range("K133").Value = 83709
$ws.Range("L133").Value = 80662.06
$ws.Range("M133").Value = -78649
$ws.Range("N133").Value = -90782.06
$ws.Range("H136").Value = 138315.33
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 138315.33
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 138315.33
$ws.Range("N136").Value = -148515.33
$ws.Range("H138").Value = 2513.4473
$ws.Range("I138").Value = 1476.9286
$ws.Range("J138").Value = 3118.0833
$ws.Range("K138").Value = 4430.7858
$ws.Range("L138").Value = 9354.249899999999
$ws.Range("M138").Value = 709.2142000000003
$ws.Range("N138").Value = -19634.2499
$ws.Range("H141").Value = 27406.027
$ws.Range("I141").Value = 29797.787
$ws.Range("J141").Value = 1096.6666
$ws.Range("K141").Value = 89393.361
$ws.Range("L141").Value = 3289.9998
$ws.Range("M141").Value = -84213.361
$ws.Range("N141").Value = -13649.9998

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 68238.516
$ws.Range("I32").Value = 49764.316
$ws.Range("J32").Value = 105186.91
$ws.Range("K32").Value = 49764.316
$ws.Range("L32").Value = 105186.91
$ws.Range("M32").Value = -49477.316
$ws.Range("N32").Value = -105760.91
$ws.Range("H74").Value = 1370.0182
$ws.Range("I74").Value = 1272.5918
$ws.Range("J74").Value = 2165.6667
$ws.Range("K74").Value = 1272.5918
$ws.Range("L74").Value = 2165.6667
$ws.Range("M74").Value = -398.5917999999999
$ws.Range("N74").Value = -3913.6667
$ws.Range("H77").Value = 1370.0182
$ws.Range("I77").Value = 1272.5918
$ws.Range("J77").Value = 2165.6667
$ws.Range("K77").Value = 6362.959
$ws.Range("L77").Value = 10828.3335
$ws.Range("M77").Value = -1994.959
$ws.Range("N77").Value = -19564.3335
$ws.Range("H122").Value = 1647.9131
$ws.Range("I122").Value = 1518.2727
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 4554.8181
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -2104.8181
$ws.Range("N122").Value = -18400
$ws.Range("H132").Value = 2076.3416
$ws.Range("I132").Value = 1137.3235
$ws.Range("J132").Value = 6637.2856
$ws.Range("K132").Value = 3411.9705
$ws.Range("L132").Value = 19911.8568
$ws.Range("M132").Value = -881.9704999999999
$ws.Range("N132").Value = -24971.8568

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 37499
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 37499
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 37499
$ws.Range("N81").Value = -39621
$ws.Range("H84").Value = 37499
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 37499
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 112497
$ws.Range("N84").Value = -123105
$ws.Range("H107").Value = 1139.6333
$ws.Range("I107").Value = 1201.56
$ws.Range("J107").Value = 830
$ws.Range("K107").Value = 1201.56
$ws.Range("L107").Value = 830
$ws.Range("M107").Value = 718.4400000000001
$ws.Range("N107").Value = -4670
$ws.Range("H134").Value = 1412.8182
$ws.Range("I134").Value = 1381.3871
$ws.Range("J134").Value = 1900
$ws.Range("K134").Value = 4144.1613
$ws.Range("L134").Value = 5700
$ws.Range("M134").Value = -1609.1613
$ws.Range("N134").Value = -10770

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 14750
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 14750
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 14750
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -15220
$ws.Range("H95").Value = 12499.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 12499.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 12499.5
$ws.Range("N95").Value = -17991.5
$ws.Range("H110").Value = 50000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 50000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
$ws.Range("H132").Value = 1794.3334
$ws.Range("I132").Value = 1723.4706
$ws.Range("J132").Value = 2999
$ws.Range("K132").Value = 5170.4118
$ws.Range("L132").Value = 8997
$ws.Range("M132").Value = -2640.4118
$ws.Range("N132").Value = -14057

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4410.8887
$ws.Range("I39").Value = 2250
$ws.Range("J39").Value = 5028.2856
$ws.Range("K39").Value = 6750
$ws.Range("L39").Value = 15084.8568
$ws.Range("M39").Value = -6456
$ws.Range("N39").Value = -15672.8568
$ws.Range("H44").Value = 215
$ws.Range("I44").Value = 172.5
$ws.Range("J44").Value = 300
$ws.Range("K44").Value = 517.5
$ws.Range("L44").Value = 900
$ws.Range("M44").Value = -119.5
$ws.Range("N44").Value = -1696
$ws.Range("H131").Value = 26267.166
$ws.Range("I131").Value = 912.6667
$ws.Range("J131").Value = 51621.668
$ws.Range("K131").Value = 2738.0001
$ws.Range("L131").Value = 154865.004
$ws.Range("M131").Value = 2301.9999
$ws.Range("N131").Value = -164945.004

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2500
$ws.Range("I31").Value = 2500
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2500
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2208
$ws.Range("H37").Value = 2500
$ws.Range("I37").Value = 2500
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2500
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2223
$ws.Range("H41").Value = 30000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 30000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 30000
$ws.Range("N41").Value = -30710
$ws.Range("H70").Value = 5824.2
$ws.Range("I70").Value = 5804.3335
$ws.Range("J70").Value = 5832.7144
$ws.Range("K70").Value = 5804.3335
$ws.Range("L70").Value = 5832.7144
$ws.Range("M70").Value = -5534.3335
$ws.Range("N70").Value = -6372.7144
$ws.Range("H73").Value = 5824.2
$ws.Range("I73").Value = 5804.3335
$ws.Range("J73").Value = 5832.7144
$ws.Range("K73").Value = 5804.3335
$ws.Range("L73").Value = 5832.7144
$ws.Range("M73").Value = -4868.3335
$ws.Range("N73").Value = -7704.7144
$ws.Range("H97").Value = 56701.25
$ws.Range("I97").Value = 110193.164
$ws.Range("J97").Value = 3209.3333
$ws.Range("K97").Value = 110193.164
$ws.Range("L97").Value = 3209.3333
$ws.Range("M97").Value = -109697.164
$ws.Range("N97").Value = -4201.3333
$ws.Range("H122").Value = 3164.125
$ws.Range("I122").Value = 3163.862
$ws.Range("J122").Value = 3166.6667
$ws.Range("K122").Value = 9491.585999999999
$ws.Range("L122").Value = 9500.000100000001
$ws.Range("M122").Value = -7041.585999999999
$ws.Range("N122").Value = -14400.0001
$ws.Range("H126").Value = 6001.6113
$ws.Range("I126").Value = 7982.8887
$ws.Range("J126").Value = 4020.3333
$ws.Range("K126").Value = 23948.6661
$ws.Range("L126").Value = 12060.9999
$ws.Range("M126").Value = -21478.6661
$ws.Range("N126").Value = -17000.9999
$ws.Range("H132").Value = 1982.4839
$ws.Range("I132").Value = 1834.32
$ws.Range("J132").Value = 2599.8333
$ws.Range("K132").Value = 5502.96
$ws.Range("L132").Value = 7799.499899999999
$ws.Range("M132").Value = -2972.96
$ws.Range("N132").Value = -12859.4999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 33432.11
$ws.Range("I7").Value = 67372.5
$ws.Range("J7").Value = 6279.8
$ws.Range("K7").Value = 67372.5
$ws.Range("L7").Value = 6279.8
$ws.Range("M7").Value = -67260.5
$ws.Range("N7").Value = -6503.8
$ws.Range("H32").Value = 2806.6667
$ws.Range("I32").Value = 2105
$ws.Range("J32").Value = 6315
$ws.Range("K32").Value = 2105
$ws.Range("L32").Value = 6315
$ws.Range("M32").Value = -1788
$ws.Range("N32").Value = -6949
$ws.Range("H40").Value = 4406.4287
$ws.Range("I40").Value = 3808.3333
$ws.Range("J40").Value = 7995
$ws.Range("K40").Value = 3808.3333
$ws.Range("L40").Value = 7995
$ws.Range("M40").Value = -3672.3333
$ws.Range("N40").Value = -8267
$ws.Range("H55").Value = 771.73334
$ws.Range("I55").Value = 651.25
$ws.Range("J55").Value = 909.4286
$ws.Range("K55").Value = 651.25
$ws.Range("L55").Value = 909.4286
$ws.Range("M55").Value = -478.25
$ws.Range("N55").Value = -1255.4286
$ws.Range("H61").Value = 1275.2941
$ws.Range("I61").Value = 1166.0667
$ws.Range("J61").Value = 2094.5
$ws.Range("K61").Value = 1166.0667
$ws.Range("L61").Value = 2094.5
$ws.Range("M61").Value = -964.0667000000001
$ws.Range("N61").Value = -2498.5
$ws.Range("H113").Value = 1275.2941
$ws.Range("I113").Value = 1166.0667
$ws.Range("J113").Value = 2094.5
$ws.Range("K113").Value = 1166.0667
$ws.Range("L113").Value = 2094.5
$ws.Range("M113").Value = 1003.9333
$ws.Range("N113").Value = -6434.5
$ws.Range("H126").Value = 33432.11
$ws.Range("I126").Value = 67372.5
$ws.Range("J126").Value = 6279.8
$ws.Range("K126").Value = 202117.5
$ws.Range("L126").Value = 18839.4
$ws.Range("M126").Value = -199647.5
$ws.Range("N126").Value = -23779.4
$ws.Range("H132").Value = 3422.2903
$ws.Range("I132").Value = 2528.3809
$ws.Range("J132").Value = 5299.5
$ws.Range("K132").Value = 7585.1427
$ws.Range("L132").Value = 15898.5
$ws.Range("M132").Value = -5055.1427
$ws.Range("N132").Value = -20958.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 62500950
$ws.Range("I107").Value = 1086.1428
$ws.Range("J107").Value = 500000000
$ws.Range("K107").Value = 3258.4284
$ws.Range("L107").Value = 1500000000
$ws.Range("M107").Value = -1338.4284
$ws.Range("N107").Value = -1500003840
$ws.Range("H132").Value = 4325.5
$ws.Range("I132").Value = 4424.7026
$ws.Range("J132").Value = 655
$ws.Range("K132").Value = 13274.1078
$ws.Range("L132").Value = 1965
$ws.Range("M132").Value = -10744.1078
$ws.Range("N132").Value = -7025

Write-Output "Applied scheduled runner price updates to all sheets."